# Issue #7 Make browse playlist mode visible
# Adds a new row (row 8) to the "Issues" log sheet recording the new issue.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")
$ws.Activate()

# New issue entry: # / Status / Name / General Settings
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "DONE"
$ws.Range("D8").Value = "Make browse playlist mode visibel"
$ws.Range("E8").Value = "Tidy UI"

# Match the row height used by the other wrapped-text issue rows.
$ws.Rows(8).RowHeight = 29

# Leave the selection where Excel would land after typing into E8.
$ws.Range("E8").Select()
